$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1351.2858
$ws.Range("I28").Value = 1077.3334
$ws.Range("J28").Value = 2995
$ws.Range("K28").Value = 1077.3334
$ws.Range("L28").Value = 2995
$ws.Range("M28").Value = -592.3334
$ws.Range("N28").Value = -3965

$ws.Range("H33").Value = 100.25
$ws.Range("I33").Value = 103.6
$ws.Range("K33").Value = 103.6
$ws.Range("M33").Value = 125.4

$ws.Range("H62").Value = 833
$ws.Range("I62").Value = 833
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 833
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -209
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 833
$ws.Range("I65").Value = 833
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 4165
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -1045
$ws.Range("N65").ClearContents()

$ws.Range("H98").Value = 1649.6428
$ws.Range("I98").Value = 1045.5714
$ws.Range("J98").Value = 2253.7144
$ws.Range("K98").Value = 1045.5714
$ws.Range("L98").Value = 2253.7144
$ws.Range("M98").Value = 452.4286
$ws.Range("N98").Value = -5249.7144

$ws.Range("H100").Value = 7024.3335
$ws.Range("I100").Value = 2536.5
$ws.Range("K100").Value = 2536.5
$ws.Range("M100").Value = -1995.5

$ws.Range("H111").Value = 2949.5
$ws.Range("I111").Value = 2949.5
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 8848.5
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = -5781.5
$ws.Range("N111").ClearContents()

$ws.Range("H122").Value = 1649.6428
$ws.Range("I122").Value = 1045.5714
$ws.Range("J122").Value = 2253.7144
$ws.Range("K122").Value = 3136.7142
$ws.Range("L122").Value = 6761.1432
$ws.Range("M122").Value = -686.7142000000003
$ws.Range("N122").Value = -11661.1432

$ws.Range("H135").Value = 890.0909
$ws.Range("I135").Value = 837.2381
$ws.Range("K135").Value = 7535.142900000001
$ws.Range("M135").Value = -5000.142900000001

$ws.Range("H141").Value = 3118.318
$ws.Range("I141").Value = 2705.6843
$ws.Range("K141").Value = 8117.0529
$ws.Range("M141").Value = -2937.0529

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4847.054
$ws.Range("I32").Value = 4680.6177
$ws.Range("K32").Value = 4680.6177
$ws.Range("M32").Value = -4393.6177

$ws.Range("H74").Value = 2687.2163
$ws.Range("I74").Value = 2376.2812
$ws.Range("K74").Value = 2376.2812
$ws.Range("M74").Value = -1502.2812

$ws.Range("H77").Value = 2687.2163
$ws.Range("I77").Value = 2376.2812
$ws.Range("K77").Value = 11881.406
$ws.Range("M77").Value = -7513.405999999999

$ws.Range("H110").Value = 2000
$ws.Range("I110").Value = 2000
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 2000
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 45
$ws.Range("N110").ClearContents()

$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1999.1666
$ws.Range("I99").Value = 1999.8
$ws.Range("J99").Value = 1996
$ws.Range("K99").Value = 1999.8
$ws.Range("L99").Value = 1996
$ws.Range("M99").Value = -501.8
$ws.Range("N99").Value = -4992

$ws.Range("H134").Value = 3308
$ws.Range("I134").Value = 3308
$ws.Range("K134").Value = 9924
$ws.Range("M134").Value = -7389

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("N12").ClearContents()

$ws.Range("I31").Value = 4999
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 4999
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -4704
$ws.Range("N31").ClearContents()

$ws.Range("I34").Value = 4999
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 4999
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -4797
$ws.Range("N34").ClearContents()

$ws.Range("H99").Value = 4184.5
$ws.Range("I99").Value = 4184.5
$ws.Range("K99").Value = 4184.5
$ws.Range("M99").Value = -2686.5

$ws.Range("H126").Value = 4184.5
$ws.Range("I126").Value = 4184.5
$ws.Range("K126").Value = 12553.5
$ws.Range("M126").Value = -10083.5

$ws.Range("H132").Value = 3323.1667
$ws.Range("I132").Value = 2986
$ws.Range("K132").Value = 8958
$ws.Range("M132").Value = -6428

$ws.Range("H134").Value = 1768.3077
$ws.Range("I134").Value = 1874
$ws.Range("K134").Value = 5622
$ws.Range("M134").Value = -3087

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 92.07692
$ws.Range("I17").Value = 92.07692
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 276.23076
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -107.23076
$ws.Range("N17").ClearContents()

$ws.Range("H55").Value = 2660
$ws.Range("J55").Value = 3175
$ws.Range("L55").Value = 9525
$ws.Range("N55").Value = -9879

$ws.Range("H109").Value = 4229
$ws.Range("I109").Value = 988.5
$ws.Range("J109").Value = 4818.1816
$ws.Range("K109").Value = 2965.5
$ws.Range("L109").Value = 14454.5448
$ws.Range("M109").Value = -1925.5
$ws.Range("N109").Value = -16534.5448

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2854.25
$ws.Range("I102").Value = 977.5
$ws.Range("J102").Value = 4731
$ws.Range("K102").Value = 977.5
$ws.Range("L102").Value = 4731
$ws.Range("M102").Value = 644.5
$ws.Range("N102").Value = -7975

$ws.Range("H122").Value = 1499.5
$ws.Range("I122").Value = 1499.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4498.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2048.5
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4500
$ws.Range("I40").Value = 4500
$ws.Range("K40").Value = 4500
$ws.Range("M40").Value = -4364

$ws.Range("H122").Value = 3003.2
$ws.Range("I122").Value = 3003.2
$ws.Range("K122").Value = 9009.599999999999
$ws.Range("M122").Value = -6559.599999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1648.381
$ws.Range("I122").Value = 1675.25
$ws.Range("K122").Value = 5025.75
$ws.Range("M122").Value = -2575.75

$ws.Range("H126").Value = 2278.1428
$ws.Range("I126").Value = 2114.75
$ws.Range("K126").Value = 6344.25
$ws.Range("M126").Value = -3874.25

$ws.Range("H132").Value = 1832.7241
$ws.Range("I132").Value = 1216.6744
$ws.Range("J132").Value = 3598.7334
$ws.Range("K132").Value = 3650.023200000001
$ws.Range("L132").Value = 10796.2002
$ws.Range("M132").Value = -1120.023200000001
$ws.Range("N132").Value = -15856.2002

$ws.Range("H136").Value = 1494.4849
$ws.Range("I136").Value = 1153.619
$ws.Range("K136").Value = 3460.857
$ws.Range("M136").Value = -910.857

Write-Host "Applied all updates"